$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-formatted Price cells to remain text (matching original inline-string cells)
$textCells = @("D5", "D6", "D10", "D11", "D12", "D15", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D35", "D38", "D39", "D42", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin, link, price and volume values
$ws.Range("D2").Value = '45.838.83'
$ws.Range("E2").Value = '  +7.76%  '
$ws.Range("D3").Value = '2.428.90'
$ws.Range("E3").Value = '  +6.32%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '114.96'
$ws.Range("E5").Value = '  +11.77%  '
$ws.Range("D6").Value = '319.17'
$ws.Range("E6").Value = '  +2.68%  '
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +5.25%  '
$ws.Range("D10").Value = '43.11'
$ws.Range("E10").Value = '  +11.11%  '
$ws.Range("D11").Value = '0.0944'
$ws.Range("E11").Value = '  +5.04%  '
$ws.Range("D12").Value = '8.75'
$ws.Range("E12").Value = '  +6.70%  '
$ws.Range("E13").Value = '  +4.95%  '
$ws.Range("E14").Value = '  +2.31%  '
$ws.Range("D15").Value = '16.04'
$ws.Range("E15").Value = '  +5.30%  '
$ws.Range("D16").Value = '2.789.06'
$ws.Range("E16").Value = '  +6.08%  '
$ws.Range("D17").Value = '2.429.54'
$ws.Range("E17").Value = '  +6.40%  '
$ws.Range("D18").Value = '45.753.49'
$ws.Range("E18").Value = '  +7.20%  '
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").Value = '  +5.06%  '
$ws.Range("E20").Value = '  +4.55%  '
$ws.Range("D21").Value = '13.45'
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '75.32'
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").Value = '3.54'
$ws.Range("E23").Value = '  +4.78%  '
$ws.Range("D24").Value = '269.42'
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").Value = '2.38'
$ws.Range("E25").Value = '  +7.41%  '
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = '7.71'
$ws.Range("E27").Value = '  +8.00%  '
$ws.Range("D28").Value = '11.39'
$ws.Range("E28").Value = '  +5.70%  '
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("D30").Value = '39.31'
$ws.Range("E30").Value = '  +10.38%  '
$ws.Range("D31").Value = '23.13'
$ws.Range("E31").Value = '  +3.53%  '
$ws.Range("D32").Value = '0.0970'
$ws.Range("E32").Value = '  +14.12%  '
$ws.Range("D33").Value = '173.13'
$ws.Range("E33").Value = '  +5.29%  '
$ws.Range("E34").Value = '  +17.00%  '
$ws.Range("D35").Value = '0.122'
$ws.Range("E35").Value = '  +9.72%  '
$ws.Range("E36").Value = '  +2.03%  '
$ws.Range("E37").Value = '  +10.85%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '4.16'
$ws.Range("E38").Value = '  +15.34%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '3.12'
$ws.Range("E39").Value = '  +12.54%  '
$ws.Range("E40").Value = '  +6.22%  '
$ws.Range("E41").Value = '  +16.69%  '
$ws.Range("D42").Value = '102.23'
$ws.Range("E42").Value = '  -4.72%  '
$ws.Range("E43").Value = '  +6.21%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").Value = '13.51'
$ws.Range("E44").Value = '  +12.23%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").Value = '71.91'
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").Value = '118.51'
$ws.Range("E47").Value = '  +7.50%  '
$ws.Range("E48").Value = '  +14.13%  '
$ws.Range("D49").Value = '1.65'
$ws.Range("E49").Value = '  +17.54%  '
$ws.Range("D50").Value = '9.44'
$ws.Range("E50").Value = '  +9.36%  '
$ws.Range("D51").Value = '79.23'
$ws.Range("E51").Value = '  +3.62%  '
